$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.010536130669065
$ws.Range("D2").Value = 1.013033400758866
$ws.Range("E2").Value = 1.01278048318499
$ws.Range("F2").Value = 1.008778288730085
$ws.Range("J2").Value = 1.015789810205607
$ws.Range("K2").Value = 1.015895768886358
$ws.Range("L2").Value = 1.015643611777393
$ws.Range("M2").Value = 1.011653506336786
$ws.Range("N2").Value = 1.017232348248819

$ws.Range("C3").Value = 1.011961825587722
$ws.Range("D3").Value = 1.014294552009238
$ws.Range("E3").Value = 1.014005743137137
$ws.Range("F3").Value = 1.010859548778239
$ws.Range("J3").Value = 1.016845823810705
$ws.Range("K3").Value = 1.01696028939188
$ws.Range("L3").Value = 1.016672285099133
$ws.Range("M3").Value = 1.0135348878805
$ws.Range("N3").Value = 1.018289861514363

$ws.Range("C4").Value = 1.012883059000038
$ws.Range("D4").Value = 1.015109664922088
$ws.Range("E4").Value = 1.014797684216438
$ws.Range("F4").Value = 1.012204372999632
$ws.Range("J4").Value = 1.017527519092688
$ws.Range("K4").Value = 1.017647633016015
$ws.Range("L4").Value = 1.017336477009675
$ws.Range("M4").Value = 1.014750042297969
$ws.Range("N4").Value = 1.018972524881837

$ws.Range("C5").Value = 1.013270046004913
$ws.Range("D5").Value = 1.015452120874031
$ws.Range("E5").Value = 1.015130410833671
$ws.Range("F5").Value = 1.012769303319354
$ws.Range("J5").Value = 1.01781372325303
$ws.Range("K5").Value = 1.017936245941869
$ws.Range("L5").Value = 1.017615367096649
$ws.Range("M5").Value = 1.015260377364592
$ws.Range("N5").Value = 1.019259135484915

$ws.Range("C6").Value = 1.013335005507032
$ws.Range("D6").Value = 1.015509608171699
$ws.Range("E6").Value = 1.0151862652314
$ws.Range("F6").Value = 1.012864132710352
$ws.Range("J6").Value = 1.017861756060277
$ws.Range("K6").Value = 1.017984685190564
$ws.Range("L6").Value = 1.017662174422573
$ws.Range("M6").Value = 1.015346035106192
$ws.Range("N6").Value = 1.019307236504258

$ws.Range("C7").Value = 1.012888231105818
$ws.Range("D7").Value = 1.015114241681824
$ws.Range("E7").Value = 1.014802130925262
$ws.Range("F7").Value = 1.01221192330768
$ws.Range("J7").Value = 1.01753134485412
$ws.Range("K7").Value = 1.017651490829374
$ws.Range("L7").Value = 1.017340204866563
$ws.Range("M7").Value = 1.014756863426336
$ws.Range("N7").Value = 1.018976356076289

$ws.Range("C8").Value = 1.011018222445679
$ws.Range("D8").Value = 1.013459810595595
$ws.Range("E8").Value = 1.013194752331359
$ws.Range("F8").Value = 1.009482063001938
$ws.Range("J8").Value = 1.01614703342633
$ws.Range("K8").Value = 1.016255837141191
$ws.Range("L8").Value = 1.015991556027772
$ws.Range("M8").Value = 1.012289798984072
$ws.Range("N8").Value = 1.017590078767489

$ws.Range("C9").Value = 1.007712793635553
$ws.Range("D9").Value = 1.010537022317888
$ws.Range("E9").Value = 1.010355294634129
$ws.Range("F9").Value = 1.004656380899225
$ws.Range("J9").Value = 1.013695035356259
$ws.Range("K9").Value = 1.013784964082866
$ws.Range("L9").Value = 1.013603862642245
$ws.Range("M9").Value = 1.007924714231784
$ws.Range("N9").Value = 1.015134598578911

$ws.Range("C10").Value = 1.005501757834982
$ws.Range("D10").Value = 1.008583054819447
$ws.Range("E10").Value = 1.008457187089967
$ws.Range("F10").Value = 1.001427778714665
$ws.Range("J10").Value = 1.012051451805072
$ws.Range("K10").Value = 1.012129558622804
$ws.Range("L10").Value = 1.012004168150265
$ws.Range("M10").Value = 1.005001624673358
$ws.Range("N10").Value = 1.01348868095056

$ws.Range("C11").Value = 1.004542476014394
$ws.Range("D11").Value = 1.007735580330921
$ws.Range("E11").Value = 1.007633978794078
$ws.Range("F11").Value = 1.000026771836208
$ws.Range("J11").Value = 1.011337556294104
$ws.Range("K11").Value = 1.011410727809409
$ws.Range("L11").Value = 1.011309526802741
$ws.Range("M11").Value = 1.003732565962925
$ws.Range("N11").Value = 1.012773771626081

$ws.Range("C12").Value = 1.004185861881782
$ws.Range("D12").Value = 1.007420572992825
$ws.Range("E12").Value = 1.007327997268281
$ws.Range("F12").Value = 0.9995059019300301
$ws.Range("J12").Value = 1.011072043372367
$ws.Range("K12").Value = 1.011143409557441
$ws.Range("L12").Value = 1.011051203965471
$ws.Range("M12").Value = 1.003260658509234
$ws.Range("N12").Value = 1.01250788164554

$ws.Range("C13").Value = 1.004262370302683
$ws.Range("D13").Value = 1.007488153114552
$ws.Range("E13").Value = 1.00739364075974
$ws.Range("F13").Value = 0.9996176520906355
$ws.Range("J13").Value = 1.011129012323739
$ws.Range("K13").Value = 1.011200764493194
$ws.Range("L13").Value = 1.011106628878238
$ws.Range("M13").Value = 1.003361908237566
$ws.Range("N13").Value = 1.012564931499357

$ws.Range("C14").Value = 1.004513004240034
$ws.Range("D14").Value = 1.007709546212349
$ws.Range("E14").Value = 1.007608690498916
$ws.Range("F14").Value = 0.9999837263603598
$ws.Range("J14").Value = 1.011315615901845
$ws.Range("K14").Value = 1.011388637618081
$ws.Range("L14").Value = 1.011288179961228
$ws.Range("M14").Value = 1.003693568724689
$ws.Range("N14").Value = 1.012751800075948

$ws.Range("C15").Value = 1.004667388790181
$ws.Range("D15").Value = 1.007845924728814
$ws.Range("E15").Value = 1.007741162302809
$ws.Range("F15").Value = 1.000209213362901
$ws.Range("J15").Value = 1.011430543189605
$ws.Range("K15").Value = 1.011504350815392
$ws.Range("L15").Value = 1.011399999326536
$ws.Range("M15").Value = 1.00389784583789
$ws.Range("N15").Value = 1.012866890573639

$ws.Range("C16").Value = 1.005565381531935
$ws.Range("D16").Value = 1.008639268849452
$ws.Range("E16").Value = 1.008511792346673
$ws.Range("F16").Value = 1.001520693768317
$ws.Range("J16").Value = 1.012098783404229
$ws.Range("K16").Value = 1.012177221654074
$ws.Range("L16").Value = 1.012050227237641
$ws.Range("M16").Value = 1.005085775770965
$ws.Range("N16").Value = 1.013536079766016

$ws.Range("C17").Value = 1.006128155681409
$ws.Range("D17").Value = 1.009136533771051
$ws.Range("E17").Value = 1.008994830672333
$ws.Range("F17").Value = 1.002342531398444
$ws.Range("J17").Value = 1.012517354827538
$ws.Range("K17").Value = 1.012598747105328
$ws.Range("L17").Value = 1.012457567408557
$ws.Range("M17").Value = 1.005830023775419
$ws.Range("N17").Value = 1.013955245608757

$ws.Range("C18").Value = 1.006456230793752
$ws.Range("D18").Value = 1.009426446076374
$ws.Range("E18").Value = 1.009276452382998
$ws.Range("F18").Value = 1.002821607394254
$ws.Range("J18").Value = 1.012761287471614
$ws.Range("K18").Value = 1.012844420486503
$ws.Range("L18").Value = 1.012694973119458
$ws.Range("M18").Value = 1.006263810019691
$ws.Range("N18").Value = 1.014199524665167

$ws.Range("C19").Value = 1.00656806545137
$ws.Range("D19").Value = 1.009525276213477
$ws.Range("E19").Value = 1.009372456905695
$ws.Range("F19").Value = 1.002984911764681
$ws.Range("J19").Value = 1.01284442631006
$ws.Range("K19").Value = 1.012928155865902
$ws.Range("L19").Value = 1.012775890503791
$ws.Range("M19").Value = 1.006411666084644
$ws.Range("N19").Value = 1.0142827815703

$ws.Range("C20").Value = 1.006067794192643
$ws.Range("D20").Value = 1.009083195861317
$ws.Range("E20").Value = 1.008943018358545
$ws.Range("F20").Value = 1.00225438595068
$ws.Range("J20").Value = 1.012472468162121
$ws.Range("K20").Value = 1.012553541667681
$ws.Range("L20").Value = 1.012413883255999
$ws.Range("M20").Value = 1.00575020636113
$ws.Range("N20").Value = 1.013910295199127

$ws.Range("C21").Value = 1.004439206996361
$ws.Range("D21").Value = 1.007644357535612
$ws.Range("E21").Value = 1.007545369440663
$ws.Range("F21").Value = 0.9998759398687239
$ws.Range("J21").Value = 1.011260675257756
$ws.Range("K21").Value = 1.011333322352852
$ws.Range("L21").Value = 1.011234726090134
$ws.Range("M21").Value = 1.00359591756811
$ws.Range("N21").Value = 1.012696781409842

$ws.Range("C22").Value = 1.003413542493865
$ws.Range("D22").Value = 1.006738439806417
$ws.Range("E22").Value = 1.00666542034573
$ws.Range("F22").Value = 0.9983777650426463
$ws.Range("J22").Value = 1.01049679989066
$ws.Range("K22").Value = 1.010564310295488
$ws.Range("L22").Value = 1.010491591642553
$ws.Range("M22").Value = 1.002238397112824
$ws.Range("N22").Value = 1.011931821252107

$ws.Range("C23").Value = 1.003957431960013
$ws.Range("D23").Value = 1.007218806407243
$ws.Range("E23").Value = 1.007132013570319
$ws.Range("F23").Value = 0.9991722443141443
$ws.Range("J23").Value = 1.010901934346577
$ws.Range("K23").Value = 1.010972152426064
$ws.Range("L23").Value = 1.010885709715521
$ws.Range("M23").Value = 1.00295833883997
$ws.Range("N23").Value = 1.012337531045421

$ws.Range("C24").Value = 1.006095069512178
$ws.Range("D24").Value = 1.009107297380737
$ws.Range("E24").Value = 1.008966430503548
$ws.Range("F24").Value = 1.002294215973764
$ws.Range("J24").Value = 1.012492751172289
$ws.Range("K24").Value = 1.012573968662617
$ws.Range("L24").Value = 1.012433622828032
$ws.Range("M24").Value = 1.00578627340651
$ws.Range("N24").Value = 1.013930607013495

$ws.Range("C25").Value = 1.008568593910982
$ws.Range("D25").Value = 1.011293561302735
$ws.Range("E25").Value = 1.011090238544764
$ws.Range("F25").Value = 1.005905871714893
$ws.Range("J25").Value = 1.014330477726549
$ws.Range("K25").Value = 1.014425152924449
$ws.Range("L25").Value = 1.014222503702876
$ws.Range("M25").Value = 1.015770943350252

